# Auto-generated edit script applying the scheduled-runner market-data update
# (Excalibur_Profits workbook) per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 5191.409
$ws.Range("I100").Value = 2069
$ws.Range("J100").Value = 7793.4165
$ws.Range("K100").Value = 2069
$ws.Range("L100").Value = 7793.4165
$ws.Range("M100").Value = -1528
$ws.Range("N100").Value = -8875.416499999999
$ws.Range("H137").Value = 639028.0600000001
$ws.Range("J137").Value = 1084662.9
$ws.Range("L137").Value = 3253988.7
$ws.Range("N137").Value = -3259088.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H54").Value = 40000
$ws.Range("J54").Value = 40000
$ws.Range("L54").Value = 40000
$ws.Range("N54").Value = -41538
$ws.Range("H61").Value = 2779187.2
$ws.Range("I61").Value = 3031659
$ws.Range("K61").Value = 3031659
$ws.Range("M61").Value = -3031447
$ws.Range("H74").Value = 2528.6155
$ws.Range("I74").Value = 987.9
$ws.Range("K74").Value = 987.9
$ws.Range("M74").Value = -113.9
$ws.Range("H77").Value = 2528.6155
$ws.Range("I77").Value = 987.9
$ws.Range("K77").Value = 4939.5
$ws.Range("M77").Value = -571.5
$ws.Range("H110").Value = 2624.75
$ws.Range("J110").Value = 3150
$ws.Range("L110").Value = 3150
$ws.Range("N110").Value = -7240
$ws.Range("H132").Value = 742857.0600000001
$ws.Range("I132").Value = 802145.6
$ws.Range("K132").Value = 2406436.8
$ws.Range("M132").Value = -2403906.8
$ws.Range("H136").Value = 2779187.2
$ws.Range("I136").Value = 3031659
$ws.Range("K136").Value = 9094977
$ws.Range("M136").Value = -9092427

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 9999
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").Value = $null
$ws.Range("H105").Value = 3531.889
$ws.Range("I105").Value = 3531.889
$ws.Range("K105").Value = 3531.889
$ws.Range("M105").Value = -1784.889
$ws.Range("H134").Value = 491069.9
$ws.Range("I134").Value = 450709.4
$ws.Range("K134").Value = 1352128.2
$ws.Range("M134").Value = -1349593.2
$ws.Range("H137").Value = 112192.5
$ws.Range("J137").Value = 112192.5
$ws.Range("L137").Value = 112192.5
$ws.Range("N137").Value = -122392.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17764.207
$ws.Range("I31").Value = 6974.647
$ws.Range("J31").Value = 33049.418
$ws.Range("K31").Value = 6974.647
$ws.Range("L31").Value = 33049.418
$ws.Range("M31").Value = -6679.647
$ws.Range("N31").Value = -33639.418
$ws.Range("H34").Value = 17764.207
$ws.Range("I34").Value = 6974.647
$ws.Range("J34").Value = 33049.418
$ws.Range("K34").Value = 6974.647
$ws.Range("L34").Value = 33049.418
$ws.Range("M34").Value = -6772.647
$ws.Range("N34").Value = -33453.418
$ws.Range("H132").Value = 11906.588
$ws.Range("I132").Value = 3122.25
$ws.Range("J132").Value = 19714.889
$ws.Range("K132").Value = 9366.75
$ws.Range("L132").Value = 59144.667
$ws.Range("M132").Value = -6836.75
$ws.Range("N132").Value = -64204.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 11096572
$ws.Range("I7").Value = 18160936
$ws.Range("K7").Value = 54482808
$ws.Range("M7").Value = -54482696
$ws.Range("H23").Value = 308.64285
$ws.Range("I23").Value = 284.4
$ws.Range("K23").Value = 853.1999999999999
$ws.Range("M23").Value = -618.1999999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 67542.5
$ws.Range("I62").Value = 55000
$ws.Range("K62").Value = 55000
$ws.Range("M62").Value = -54314
$ws.Range("H65").Value = 67542.5
$ws.Range("I65").Value = 55000
$ws.Range("K65").Value = 165000
$ws.Range("M65").Value = -161568
$ws.Range("H70").Value = 5052.273
$ws.Range("I70").Value = 5045
$ws.Range("J70").Value = 5056.4287
$ws.Range("K70").Value = 5045
$ws.Range("L70").Value = 5056.4287
$ws.Range("M70").Value = -4775
$ws.Range("N70").Value = -5596.4287
$ws.Range("H73").Value = 5052.273
$ws.Range("I73").Value = 5045
$ws.Range("J73").Value = 5056.4287
$ws.Range("K73").Value = 5045
$ws.Range("L73").Value = 5056.4287
$ws.Range("M73").Value = -4109
$ws.Range("N73").Value = -6928.4287
$ws.Range("H123").Value = 74998.664
$ws.Range("J123").Value = 74998.664
$ws.Range("L123").Value = 74998.664
$ws.Range("N123").Value = -79898.664
$ws.Range("H132").Value = 32651040
$ws.Range("I132").Value = 42172610
$ws.Range("J132").Value = 5660.4287
$ws.Range("K132").Value = 126517830
$ws.Range("L132").Value = 16981.2861
$ws.Range("M132").Value = -126515300
$ws.Range("N132").Value = -22041.2861
$ws.Range("H141").Value = 22500
$ws.Range("J141").Value = 22500
$ws.Range("L141").Value = 22500
$ws.Range("N141").Value = -32860

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 40000
$ws.Range("J54").Value = 40000
$ws.Range("L54").Value = 40000
$ws.Range("N54").Value = -41288
$ws.Range("H61").Value = 1949.75
$ws.Range("I61").Value = 1679.8
$ws.Range("K61").Value = 1679.8
$ws.Range("M61").Value = -1477.8
$ws.Range("H93").Value = 1709.25
$ws.Range("I93").Value = 1712.4667
$ws.Range("J93").Value = 1699.6
$ws.Range("K93").Value = 1712.4667
$ws.Range("L93").Value = 1699.6
$ws.Range("M93").Value = -464.4666999999999
$ws.Range("N93").Value = -4195.6
$ws.Range("H100").Value = 9907.77
$ws.Range("I100").Value = 2345.182
$ws.Range("J100").Value = 51502
$ws.Range("K100").Value = 2345.182
$ws.Range("L100").Value = 51502
$ws.Range("M100").Value = -1804.182
$ws.Range("N100").Value = -52584
$ws.Range("H113").Value = 1949.75
$ws.Range("I113").Value = 1679.8
$ws.Range("K113").Value = 1679.8
$ws.Range("M113").Value = 490.2
$ws.Range("H132").Value = 1517829.8
$ws.Range("I132").Value = 2681369
$ws.Range("J132").Value = 5228.7
$ws.Range("K132").Value = 8044107
$ws.Range("L132").Value = 15686.1
$ws.Range("M132").Value = -8041577
$ws.Range("N132").Value = -20746.1
$ws.Range("H133").Value = 51630.332
$ws.Range("J133").Value = 51630.332
$ws.Range("L133").Value = 51630.332
$ws.Range("N133").Value = -56690.332
$ws.Range("H136").Value = 4179.636
$ws.Range("I136").Value = 2997
$ws.Range("J136").Value = 6249.25
$ws.Range("K136").Value = 8991
$ws.Range("L136").Value = 18747.75
$ws.Range("M136").Value = -6441
$ws.Range("N136").Value = -23847.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5921576.5
$ws.Range("I132").Value = 8051984
$ws.Range("J132").Value = 3777
$ws.Range("K132").Value = 24155952
$ws.Range("L132").Value = 11331
$ws.Range("M132").Value = -24153422
$ws.Range("N132").Value = -16391
$ws.Range("H136").Value = 22266.291
$ws.Range("I136").Value = 30899.812
$ws.Range("J136").Value = 4999.25
$ws.Range("K136").Value = 92699.436
$ws.Range("L136").Value = 14997.75
$ws.Range("M136").Value = -90149.436
$ws.Range("N136").Value = -20097.75
